# Apply "pert mmd" edit: change the PI() coefficient to (PI()+e) and update
# the underlying naive/PERT totals on both sheets, then make the "Наивный"
# sheet the active tab.

$wb = $excel.ActiveWorkbook

$wsNaive = $wb.Worksheets.Item("Наивный")
$wsPert  = $wb.Worksheets.Item("PERT")

# ---- Sheet "Наивный" ----------------------------------------------------
# Update the "ИТОГО" (optimistic/pessimistic/base) totals rows (30-31)
$wsNaive.Range("D30").Value = 600
$wsNaive.Range("E30").Value = 1000
$wsNaive.Range("F30").Value = 700

$wsNaive.Range("D31").Value = 400
$wsNaive.Range("E31").Value = 800
$wsNaive.Range("F31").Value = 500

# Update the PI() formula to (PI()+2.71)
$wsNaive.Range("D35").Formula = "=D34*(PI()+2.71)"
$wsNaive.Range("E35").Formula = "=E34*(PI()+2.71)"
$wsNaive.Range("F35").Formula = "=F34*(PI()+2.71)"

# ---- Sheet "PERT" --------------------------------------------------------
$wsPert.Range("C30").Value = 600
$wsPert.Range("D30").Value = 1000
$wsPert.Range("E30").Value = 700

$wsPert.Range("C31").Value = 400
$wsPert.Range("D31").Value = 800
$wsPert.Range("E31").Value = 500

$wsPert.Range("E36").Formula = "=E35*(PI()+2.71)"

# ---- Shared label text ----------------------------------------------------
$wsNaive.Range("A35").Value = "ИТОГО (человек/месяц) * (π+e)"
$wsPert.Range("A36").Value = "ИТОГО (человек/месяц) * (π+e)"

# ---- Active sheet / selection -------------------------------------------
$wsNaive.Select()
$wsNaive.Range("A1:F35").Select()
$wsNaive.Range("F35").Activate()

$wsPert.Range("A1:G36").Select()
$wsPert.Range("F32").Activate()

$wsNaive.Select()
